# 'Edit profile' module updated.
# Switch the active/selected sheet from "SIGN UP" to "Sheet3" (the Edit
# profile sheet), update the edited profile's Name / DOB / image path,
# widen column B, reset row 13's custom height, and move the selection.

$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)

# Update the profile values shown on the "Edit profile" sheet.
$ws3.Cells.Item(1, 2).Value = " Thalwar"

# The dob cell is formatted as a date (with quotePrefix) - a leading
# apostrophe forces the new value to stay literal text instead of being
# reinterpreted by Excel as a real date (which would change its format).
$ws3.Cells.Item(2, 2).Value = "'22/04/2028"

$ws3.Cells.Item(13, 2).Value = "src\main\resources\sample%20(1).pdf"

# Widen column B on Sheet3.
$ws3.Columns.Item(2).ColumnWidth = 42.8

# Row 13 no longer needs the taller custom height - restore auto height.
$ws3.Rows.Item(13).AutoFit()

# Make Sheet3 the active/selected tab, with D7 as the active selected cell.
$ws3.Activate() | Out-Null
$ws3.Range("D7").Select() | Out-Null
